# Apply the "cryptos list" update (GitHub Actions data refresh).
# Only cell values in columns B-E change; row index column A and headers are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.949.45'
$ws.Range("E2").Value = '  -1.31%  '
# Row 3
$ws.Range("D3").Value = '3.371.03'
$ws.Range("E3").Value = '  -0.61%  '
# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
# Row 5
$ws.Range("D5").Value = '''572.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.89%  '
# Row 6
$ws.Range("D6").Value = '''137.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.16%  '
# Row 7
$ws.Range("E7").Value = '  -0.01%  '
# Row 8
$ws.Range("D8").Value = '3.371.99'
$ws.Range("E8").Value = '  -0.52%  '
# Row 9
$ws.Range("E9").Value = '  -1.13%  '
# Row 10
$ws.Range("D10").Value = '''7.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.91%  '
# Row 11
$ws.Range("E11").Value = '  -2.82%  '
# Row 12
$ws.Range("E12").Value = '  -1.70%  '
# Row 13
$ws.Range("D13").Value = '3.942.79'
$ws.Range("E13").Value = '  -0.72%  '
# Row 14
$ws.Range("E14").Value = '  +0.71%  '
# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000172'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.33%  '
# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.370.45'
$ws.Range("E17").Value = '  -0.60%  '
# Row 18
$ws.Range("D18").Value = '60.989.13'
$ws.Range("E18").Value = '  -1.34%  '
# Row 19
$ws.Range("D19").Value = '''13.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.75%  '
# Row 20
$ws.Range("D20").Value = '''5.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.69%  '
# Row 21
$ws.Range("D21").Value = '''9.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.69%  '
# Row 22
$ws.Range("D22").Value = '''375.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.27%  '
# Row 23
$ws.Range("D23").Value = '''0.550'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.00%  '
# Row 24
$ws.Range("D24").Value = '3.506.21'
$ws.Range("E24").Value = '  -0.69%  '
# Row 25
$ws.Range("E25").Value = '  +0.07%  '
# Row 26
$ws.Range("D26").Value = '''0.0000126'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.28%  '
# Row 27
$ws.Range("D27").Value = '''70.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.56%  '
# Row 28
$ws.Range("D28").Value = '''1.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.30%  '
# Row 29
$ws.Range("D29").Value = '''0.178'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +10.12%  '
# Row 30
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '''7.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.86%  '
# Row 31
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = '''1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.01%  '
# Row 32
$ws.Range("D32").Value = '''8.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.02%  '
# Row 33
$ws.Range("D33").Value = '''2.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.57%  '
# Row 34
$ws.Range("E34").Value = '  -0.04%  '
# Row 35
$ws.Range("D35").Value = '''23.64'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.74%  '
# Row 36
$ws.Range("D36").Value = '''5.20'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.36%  '
# Row 37
$ws.Range("D37").Value = '''6.85'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.15%  '
# Row 38
$ws.Range("E38").Value = '  -0.90%  '
# Row 39
$ws.Range("D39").Value = '''164.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.07%  '
# Row 40
$ws.Range("D40").Value = '''0.0760'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.69%  '
# Row 41
$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.16%  '
# Row 42
$ws.Range("E42").Value = '  -1.22%  '
# Row 43
$ws.Range("D43").Value = '''41.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
# Row 44
$ws.Range("D44").Value = '''1.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.73%  '
# Row 45
$ws.Range("D45").Value = '''4.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.42%  '
# Row 46
$ws.Range("E46").Value = '  -3.65%  '
# Row 47
$ws.Range("D47").Value = '''24.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.23%  '
# Row 48
$ws.Range("D48").Value = '2.452.41'
$ws.Range("E48").Value = '  +4.37%  '
# Row 49
$ws.Range("D49").Value = '''6.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.72%  '
# Row 50
$ws.Range("D50").Value = '''22.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.88%  '
# Row 51
$ws.Range("E51").Value = '  +3.58%  '
